# Insert two new price-record rows at the top of the "Vega Modelo de Temuco -
# Coliflor" weekly data block (rows 601-602), pushing the existing rows
# 601:653 down to 603:655.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("601:602").Insert()

# New row 601
$ws.Range("A601").Value = 10
$ws.Range("B601").Value = "Vega Modelo de Temuco"
$ws.Range("C601").Value = "La Araucanía"
$ws.Range("D601").Value = 45106
$ws.Range("E601").Value = 9
$ws.Range("F601").Value = 100112008
$ws.Range("G601").Value = "Coliflor"
$ws.Range("H601").Value = "Sin especificar"
$ws.Range("I601").Value = "Primera"
$ws.Range("J601").Value = 2800
$ws.Range("K601").Value = 1200
$ws.Range("L601").Value = 1200
$ws.Range("M601").Value = 1200
$ws.Range("N601").Value = "`$/unidad"
$ws.Range("O601").Value = "Región Metropolitana"
$ws.Range("P601").Value = 1200
$ws.Range("Q601").Value = 1
$ws.Range("R601").Value = "Hortaliza"

# New row 602
$ws.Range("A602").Value = 10
$ws.Range("B602").Value = "Vega Modelo de Temuco"
$ws.Range("C602").Value = "La Araucanía"
$ws.Range("D602").Value = 45106
$ws.Range("E602").Value = 9
$ws.Range("F602").Value = 100112008
$ws.Range("G602").Value = "Coliflor"
$ws.Range("H602").Value = "Sin especificar"
$ws.Range("I602").Value = "Primera"
$ws.Range("J602").Value = 2100
$ws.Range("K602").Value = 1300
$ws.Range("L602").Value = 1300
$ws.Range("M602").Value = 1300
$ws.Range("N602").Value = "`$/unidad"
$ws.Range("O602").Value = "Región de O'Higgins"
$ws.Range("P602").Value = 1300
$ws.Range("Q602").Value = 1
$ws.Range("R602").Value = "Hortaliza"

# Make sure the date-column number format matches the rest of column D.
$ws.Range("D601:D602").NumberFormat = $ws.Range("D603").NumberFormat
